# Add a new participant record (bkp028) into the WJ scores table.
# Excel was originally positioned with row 7 (bkp027) followed by bkp029;
# a row is inserted above the old row 8 so bkp028 sits in ID order and every
# row below shifts down by one (absorbing the single blank spacer row that
# used to follow the last participant, bkp120).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 8 - everything from the old row 8 down
# (through the former blank row 66) shifts down to make room.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row with the new participant's data.
$ws.Range("A8").Value = "bkp028"
$ws.Range("B8").Value = 68

# Leave the selection where the author left it when saving.
$ws.Range("A9").Select()
